$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 24 with the new task entry
$ws.Range("A24").Value = "esthétique globale"
$ws.Range("B24").Value = "Roméo(directeur artistique)"

# Match the style (centered) used by the rest of column A/B cells
$ws.Range("A24").HorizontalAlignment = -4108
$ws.Range("B24").HorizontalAlignment = -4108

# Update the active selection to B24 as in the final workbook
$ws.Range("B24").Select()
